$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the 1-based index of the paragraph whose visible text (not
# counting the trailing paragraph-mark / cell-mark control chars) equals
# $text exactly. Returns -1 if not found.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($text) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        $ptext = $p.Range.Text
        $ptext = $ptext.TrimEnd([char]13, [char]7)
        if ($ptext -eq $text) {
            return $i
        }
    }
    return -1
}

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1. The "BootStrap modal" reference link paragraph currently carries the
#    stray _GoBack bookmark at its end. Rewrite it without the bookmark,
#    keeping its existing list formatting and text untouched.
# ---------------------------------------------------------------------------
$bootstrapIdx = Find-ParagraphIndex("https://getbootstrap.com/docs/4.0/components/modal/")
if ($bootstrapIdx -eq -1) {
    throw "Could not locate the getbootstrap.com reference paragraph"
}
$bootstrapPara = $d.Paragraphs($bootstrapIdx)
$bootstrapXml = "<w:p $wNs>" + `
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='4'/></w:numPr></w:pPr>" + `
    "<w:r><w:t>https://getbootstrap.com/docs/4.0/components/modal/</w:t></w:r>" + `
    "</w:p>"
[void]$bootstrapPara.Range.InsertXML($bootstrapXml)

# ---------------------------------------------------------------------------
# 2. The "This website" list item (currently highlighted yellow) becomes the
#    CodePen reference link, loses the yellow highlight, and gains the
#    _GoBack bookmark that used to sit on the BootStrap paragraph.
# ---------------------------------------------------------------------------
$thisWebsiteIdx = Find-ParagraphIndex("This website")
if ($thisWebsiteIdx -eq -1) {
    throw "Could not locate the 'This website' paragraph"
}
$thisWebsitePara = $d.Paragraphs($thisWebsiteIdx)
$codepenXml = "<w:p $wNs>" + `
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='4'/></w:numPr></w:pPr>" + `
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" + `
    "<w:r><w:t>https://codepen.io/lucy_wheel/pen/VxYzKP</w:t></w:r>" + `
    "</w:p>"
[void]$thisWebsitePara.Range.InsertXML($codepenXml)
